$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B, C, D, E, F, H across rows 2-9 (column G unchanged)
$data = @{
    2 = @{ B = 1168.76318359375;  C = 0.8768;  D = 0.8686000108718872; E = 1.129199981689453;  F = 0.6503999829292297; H = 0.3817 }
    3 = @{ B = 1353.689453125;    C = 1.083;   D = 0.9636;             E = 2.630399942398071;  F = 0.6661999821662903; H = 1.2235 }
    4 = @{ B = 889.8314208984375; C = 1.0395;  D = 0.9536;             E = 2.359600067138672;  F = 0.7373999953269958; H = 1.1348 }
    5 = @{ B = 870.61767578125;   C = 0.9184;  D = 0.8925999999999999; E = 1.464900016784668;  F = 0.6840999722480774; H = 0.5946 }
    6 = @{ B = 1139.08642578125;  C = 0.9012;  D = 0.8973;             E = 1.151800036430359;  F = 0.7023000121116638; H = 0.6358 }
    7 = @{ B = 859.8021240234375; C = 0.8659;  D = 0.8579999804496765; E = 1.12090003490448;   F = 0.72079998254776;   H = 0.288  }
    8 = @{ B = 942.6757202148438; C = 0.8447;  D = 0.839;              E = 1.106600046157837;  F = 0.7278000116348267; H = 0.1195 }
    9 = @{ B = 7224.46630859375;  C = 0.9311;  D = 0.8848;             E = 2.630399942398071;  F = 0.6503999829292297; H = 4.3779 }
}

foreach ($row in $data.Keys) {
    $rowData = $data[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$row").Value = $rowData[$col]
    }
}
